$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics with refreshed data values
$ws.Range("C8").Value = 1441
$ws.Range("D8").Value = 221
$ws.Range("E8").Value = 1220
$ws.Range("F8").Value = 9.064807219031994
$ws.Range("G8").Value = 84.66342817487855
$ws.Range("H8").Value = 15.33657182512144
